# Generate Report for Handoff
# Re-sort the "In Translation" / "Ready for handoff" block of rows (originally
# 69122473, 86f42771, fb659db4, 0689cad6) in each of the three sheets so that
# 69122473 moves to the back of the block and becomes "Ready for handoff"
# with a fresh handoff timestamp, while 86f42771 / fb659db4 / 0689cad6 shift
# up one row each.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A5").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.md"
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"
$wsOverview.Range("D5").Value = "2016-03-22 05:03:26"

$wsOverview.Range("A6").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.md"
$wsOverview.Range("B6").Value = "In Translation"
$wsOverview.Range("C6").Value = "In Translation"
$wsOverview.Range("D6").Value = "2016-03-22 05:02:02"

$wsOverview.Range("A7").Value = "0689cad6-d415-4d7a-87de-99acd646a097.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-03-22 05:11:17"

$wsOverview.Range("A8").Value = "69122473-a34a-4d10-b38d-1fdd26acc32e.md"
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"
$wsOverview.Range("D8").Value = "2016-03-22 05:15:10"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A5").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("D5").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.zh-cn.xlf"
$wsZhCn.Range("E5").Value = "2016-03-22 05:03:18"

$wsZhCn.Range("A6").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.md"
$wsZhCn.Range("C6").Value = "In Translation"
$wsZhCn.Range("D6").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-22 05:01:47"

$wsZhCn.Range("A7").Value = "0689cad6-d415-4d7a-87de-99acd646a097.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "0689cad6-d415-4d7a-87de-99acd646a097.d630e21bfcebe14a13ae89e1e412fcd9f1e4af7e.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-22 05:11:10"

$wsZhCn.Range("A8").Value = "69122473-a34a-4d10-b38d-1fdd26acc32e.md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
$wsZhCn.Range("D8").Value = "69122473-a34a-4d10-b38d-1fdd26acc32e.a71666560ccbe46a04786a73c5382332c1b52989.zh-cn.xlf"
$wsZhCn.Range("E8").Value = "2016-03-22 05:15:01"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A5").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("D5").Value = "86f42771-06db-4da6-93e3-4a9e101966cb.3ca688be53f46e1fac1d7c40e4b6d84b3c7b8f04.de-de.xlf"
$wsDeDe.Range("E5").Value = "2016-03-22 05:03:26"

$wsDeDe.Range("A6").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.md"
$wsDeDe.Range("C6").Value = "In Translation"
$wsDeDe.Range("D6").Value = "fb659db4-70f5-4538-8936-7ccd74a12800.7a592e1fd0fc852671a53e7f3b03df8e95b02793.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-22 05:02:02"

$wsDeDe.Range("A7").Value = "0689cad6-d415-4d7a-87de-99acd646a097.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "0689cad6-d415-4d7a-87de-99acd646a097.d630e21bfcebe14a13ae89e1e412fcd9f1e4af7e.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-22 05:11:17"

$wsDeDe.Range("A8").Value = "69122473-a34a-4d10-b38d-1fdd26acc32e.md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
$wsDeDe.Range("D8").Value = "69122473-a34a-4d10-b38d-1fdd26acc32e.a71666560ccbe46a04786a73c5382332c1b52989.de-de.xlf"
$wsDeDe.Range("E8").Value = "2016-03-22 05:15:10"
